$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove Nancy Pelosi's row (was row 227: id P000197) -- she is excluded
# ("SOH") from the cosponsor analysis. Deleting the row shifts every
# subsequent record up by one (228->227, ..., 236->235) and the sheet's
# used range shrinks from BM236 to BM235.
$ws.Rows("227").Delete()

# Normalize the "stance" column's non-sponsoring label from a space to an
# underscore for every remaining data row.
$ws.Range("A210:A235").Replace("not sponsoring", "not_sponsoring")
